$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.628.52"
$ws.Range("E2").Value = "  +2.66%  "

$ws.Range("D3").Value = "1.890.21"
$ws.Range("E3").Value = "  +1.29%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.06"
$ws.Range("E5").Value = "  +2.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4930"
$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2954"
$ws.Range("E8").Value = "  +1.83%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06783"
$ws.Range("E9").Value = "  +3.19%  "

$ws.Range("D10").Value = "1.889.99"
$ws.Range("E10").Value = "  +1.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.36"
$ws.Range("E11").Value = "  +4.75%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07326"
$ws.Range("E12").Value = "  +1.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "90.58"
$ws.Range("E13").Value = "  +5.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.087"
$ws.Range("E14").Value = "  +4.68%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6789"
$ws.Range("E15").Value = "  +2.35%  "

$ws.Range("D16").Value = "30.627.69"
$ws.Range("E16").Value = "  +2.73%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008008"
$ws.Range("E17").Value = "  +2.66%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.24"
$ws.Range("E19").Value = "  +4.43%  "

$ws.Range("D20").Value = "2.126.47"
$ws.Range("E20").Value = "  +0.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9995"
$ws.Range("E21").Value = "  +0.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.860"
$ws.Range("E22").Value = "  +2.77%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "189.76"
$ws.Range("E23").Value = "  +38.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.050"
$ws.Range("E24").Value = "  +8.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.326"
$ws.Range("E25").Value = "  +3.59%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.46"
$ws.Range("E26").Value = "  +2.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.16"
$ws.Range("E27").Value = "  +13.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.926"
$ws.Range("E28").Value = "  +1.62%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.390"
$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.330"
$ws.Range("E30").Value = "  +3.67%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09026"
$ws.Range("E31").Value = "  +3.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.038"
$ws.Range("E32").Value = "  +2.45%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05185"
$ws.Range("E33").Value = "  +3.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7440"
$ws.Range("E34").Value = "  +5.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.128"
$ws.Range("E35").Value = "  +3.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.698"
$ws.Range("E36").Value = "  +1.44%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01882"
$ws.Range("E37").Value = "  +9.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.674"
$ws.Range("E38").Value = "  +0.15%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.156"
$ws.Range("E39").Value = "  -0.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9384"
$ws.Range("E40").Value = "  +1.85%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4463"
$ws.Range("E41").Value = "  +5.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.98"
$ws.Range("E42").Value = "  +4.51%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.813"
$ws.Range("E43").Value = "  +0.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").Value = "  +0.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.673"
$ws.Range("E45").Value = "  +3.91%  "

$ws.Range("E46").Value = "  +7.53%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05847"
$ws.Range("E47").Value = "  +3.76%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.3974"
$ws.Range("E48").Value = "  +6.04%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.47"
$ws.Range("E49").Value = "  +3.95%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.569"
$ws.Range("E50").Value = "  +6.37%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.407"
$ws.Range("E51").Value = "  +6.31%  "
